$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5737.706
$ws.Range("I137").Value = 1153.4166
$ws.Range("J137").Value = 16740
$ws.Range("K137").Value = 3460.2498
$ws.Range("L137").Value = 50220
$ws.Range("M137").Value = -910.2498000000001
$ws.Range("N137").Value = -55320

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1077.3334
$ws.Range("I141").Value = 886.4286
$ws.Range("K141").Value = 2659.2858
$ws.Range("M141").Value = 2520.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3854.2546
$ws.Range("I134").Value = 2594.6875
$ws.Range("J134").Value = 5606.696
$ws.Range("K134").Value = 7784.0625
$ws.Range("L134").Value = 16820.088
$ws.Range("M134").Value = -5249.0625
$ws.Range("N134").Value = -21890.088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3729.4465
$ws.Range("I31").Value = 2220.3447
$ws.Range("J31").Value = 5350.3335
$ws.Range("K31").Value = 2220.3447
$ws.Range("L31").Value = 5350.3335
$ws.Range("M31").Value = -1925.3447
$ws.Range("N31").Value = -5940.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3729.4465
$ws.Range("I34").Value = 2220.3447
$ws.Range("J34").Value = 5350.3335
$ws.Range("K34").Value = 2220.3447
$ws.Range("L34").Value = 5350.3335
$ws.Range("M34").Value = -2018.3447
$ws.Range("N34").Value = -5754.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 24500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9346.25
$ws.Range("J50").Value = 9346.25
$ws.Range("L50").Value = 9346.25
$ws.Range("N50").Value = -10596.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9282.666999999999
$ws.Range("J51").Value = 9282.666999999999
$ws.Range("L51").Value = 9282.666999999999
$ws.Range("N51").Value = -10754.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 15714.5
$ws.Range("J59").Value = 15714.5
$ws.Range("L59").Value = 15714.5
$ws.Range("N59").Value = -18004.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 8013.1665
$ws.Range("I60").Value = 7500
$ws.Range("J60").Value = 8269.75
$ws.Range("K60").Value = 7500
$ws.Range("L60").Value = 8269.75
$ws.Range("M60").Value = -6989
$ws.Range("N60").Value = -9291.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9282.666999999999
$ws.Range("J61").Value = 9282.666999999999
$ws.Range("L61").Value = 9282.666999999999
$ws.Range("N61").Value = -9978.666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17247.25
$ws.Range("J68").Value = 17247.25
$ws.Range("L68").Value = 17247.25
$ws.Range("N68").Value = -18745.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17247.25
$ws.Range("J71").Value = 17247.25
$ws.Range("L71").Value = 51741.75
$ws.Range("N71").Value = -59229.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 13777.25
$ws.Range("J74").Value = 13777.25
$ws.Range("L74").Value = 13777.25
$ws.Range("N74").Value = -15525.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 13777.25
$ws.Range("J77").Value = 13777.25
$ws.Range("L77").Value = 41331.75
$ws.Range("N77").Value = -50067.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2842.4167
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 3010.9
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 3010.9
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -6006.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2842.4167
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3010.9
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 9032.700000000001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13972.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 992243.4
$ws.Range("I2").Value = 135.44444
$ws.Range("J2").Value = 2778037.5
$ws.Range("K2").Value = 812.6666399999999
$ws.Range("L2").Value = 16668225
$ws.Range("M2").Value = -699.6666399999999
$ws.Range("N2").Value = -16668451

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 652.4286
$ws.Range("I34").Value = 126
$ws.Range("J34").Value = 1600
$ws.Range("K34").Value = 378
$ws.Range("L34").Value = 4800
$ws.Range("M34").Value = -294
$ws.Range("N34").Value = -4968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3300
$ws.Range("J46").Value = 3300
$ws.Range("L46").Value = 9900
$ws.Range("N46").Value = -10082

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3015
$ws.Range("I58").Value = 310
$ws.Range("J58").Value = 3916.6667
$ws.Range("K58").Value = 930
$ws.Range("L58").Value = 11750.0001
$ws.Range("M58").Value = -802
$ws.Range("N58").Value = -12006.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 9898
$ws.Range("J64").Value = 13333.333
$ws.Range("L64").Value = 39999.999
$ws.Range("N64").Value = -40539.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 9898
$ws.Range("J67").Value = 13333.333
$ws.Range("L67").Value = 39999.999
$ws.Range("N67").Value = -41871.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1795.1571
$ws.Range("I131").Value = 2985.1177
$ws.Range("J131").Value = 1413.4717
$ws.Range("K131").Value = 8955.3531
$ws.Range("L131").Value = 4240.4151
$ws.Range("M131").Value = -3915.3531
$ws.Range("N131").Value = -14320.4151

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4679.222
$ws.Range("I102").Value = 7594.5
$ws.Range("J102").Value = 2347
$ws.Range("K102").Value = 7594.5
$ws.Range("L102").Value = 2347
$ws.Range("M102").Value = -5972.5
$ws.Range("N102").Value = -5591

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3870.2917
$ws.Range("I126").Value = 3025.8333
$ws.Range("J126").Value = 4714.75
$ws.Range("K126").Value = 9077.499899999999
$ws.Range("L126").Value = 14144.25
$ws.Range("M126").Value = -6607.499899999999
$ws.Range("N126").Value = -19084.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2989.0286
$ws.Range("I40").Value = 6462.1
$ws.Range("J40").Value = 1599.8
$ws.Range("K40").Value = 6462.1
$ws.Range("L40").Value = 1599.8
$ws.Range("M40").Value = -6326.1
$ws.Range("N40").Value = -1871.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1921.25
$ws.Range("I46").Value = 2798
$ws.Range("J46").Value = 1629
$ws.Range("K46").Value = 2798
$ws.Range("L46").Value = 1629
$ws.Range("M46").Value = -2610
$ws.Range("N46").Value = -2005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1785.8545
$ws.Range("I132").Value = 1518.4445
$ws.Range("J132").Value = 2989.2
$ws.Range("K132").Value = 4555.333500000001
$ws.Range("L132").Value = 8967.599999999999
$ws.Range("M132").Value = -2025.333500000001
$ws.Range("N132").Value = -14027.6
